{"js": "// Highlight quantitative metrics (percentages, dollar amounts, large\n// numbers) in specific bullet/paragraph lines by splitting the run that\n// contains each metric into its own run and applying bold + a dark\n// slate color (#2C3E50) to it - matching the \"hybrid bold + color\n// highlighting\" described in the commit message.\n\nconst HIGHLIGHT_COLOR = \"2C3E50\";\n\n// Each entry identifies the target paragraph by a unique substring of its\n// *original* text, plus the ordered list of metric substrings inside that\n// paragraph which must become bold + colored runs.\nconst EDITS = [\n  {\n    contains: \"Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms\",\n    metrics: [\"23%\", \"64%\"]\n  },\n  {\n    contains: \"Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins\",\n    metrics: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"]\n  },\n  {\n    contains: \"Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"]\n  },\n  {\n    contains: \"Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database\",\n    metrics: [\"$400M\", \"$1B\"]\n  },\n  {\n    contains: \"Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"]\n  },\n  {\n    contains: \"Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"]\n  }\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Snapshot original text per-paragraph up front since later formatting\n// edits split runs but never change the visible characters, so matching\n// against a pre-loaded text index stays valid across the whole pass.\nconst paraTexts = paragraphs.items.map((p) => p.text);\nconst used = new Array(paraTexts.length).fill(false);\n\nfor (const edit of EDITS) {\n  let targetIndex = -1;\n  for (let i = 0; i < paraTexts.length; i++) {\n    if (!used[i] && paraTexts[i].indexOf(edit.contains) >= 0) {\n      targetIndex = i;\n      break;\n    }\n  }\n  if (targetIndex === -1) {\n    throw new Error(\"Could not locate target paragraph for: \" + edit.contains);\n  }\n  used[targetIndex] = true;\n\n  const paragraph = paragraphs.items[targetIndex];\n  for (const metric of edit.metrics) {\n    const found = paragraph.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length === 0) {\n      throw new Error(\"Could not locate metric '\" + metric + \"' in paragraph: \" + edit.contains);\n    }\n    const hit = found.items[0];\n    hit.font.bold = true;\n    hit.font.color = HIGHLIGHT_COLOR;\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight quantitative metrics (percentages, dollar amounts, large\n# numbers) in specific bullet/paragraph lines by splitting the run that\n# contains each metric into its own run and applying bold + a dark\n# slate color (#2C3E50) to it - matching the \"hybrid bold + color\n# highlighting\" described in the commit message.\n#\n# Word's Font.Color takes a packed BGR integer (classic VBA wdColor /\n# RGB() packing: R + G*256 + B*65536), so #2C3E50 (R=0x2C,G=0x3E,B=0x50)\n# becomes 0x2C + 0x3E*256 + 0x50*65536 = 5258796.\n\n$d = $word.ActiveDocument\n$highlightColor = 5258796\n\n$edits = @(\n    @{\n        Contains = '*Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms*'\n        Metrics  = @('23%', '64%')\n    },\n    @{\n        Contains = '*Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins*'\n        Metrics  = @('87%', '71%', '\u00b14.2%', '\u00b12.1%')\n    },\n    @{\n        Contains = '*Wrote RFP and analyzed bids from 1,200 vendors for research platform development*'\n        Metrics  = @('1,200')\n    },\n    @{\n        Contains = '*Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database*'\n        Metrics  = @('$400M', '$1B')\n    },\n    @{\n        Contains = '*Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M*'\n        Metrics  = @('73.5%', '$4.7M')\n    },\n    @{\n        Contains = '*Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%*'\n        Metrics  = @('87%', '71%')\n    }\n)\n\n$paraCount = $d.Paragraphs.Count\n$usedParagraphs = @()\n\nforeach ($edit in $edits) {\n    $targetIndex = -1\n    for ($i = 1; $i -le $paraCount; $i++) {\n        if ($usedParagraphs -contains $i) {\n            continue\n        }\n        $paraText = $d.Paragraphs.Item($i).Range.Text\n        if ($paraText -like $edit.Contains) {\n            $targetIndex = $i\n            break\n        }\n    }\n    if ($targetIndex -eq -1) {\n        Write-Output \"Could not locate target paragraph for: $($edit.Contains)\"\n        continue\n    }\n    $usedParagraphs += $targetIndex\n\n    foreach ($metric in $edit.Metrics) {\n        $searchRange = $d.Paragraphs.Item($targetIndex).Range.Duplicate\n        $found = $searchRange.Find.Execute($metric)\n        if (-not $found) {\n            Write-Output \"Could not locate metric '$metric' in paragraph $targetIndex\"\n            continue\n        }\n        $searchRange.Font.Bold = 1\n        $searchRange.Font.Color = $highlightColor\n    }\n}\n"}
